# tidyReport.xlsx update:
#  1. Query start/end optional (labels for 開始:/結束: stay the same text,
#     they just naturally get re-indexed in sharedStrings.xml).
#  2. Update tidy report header row: insert 訂單/品名 columns, rename
#     流動卡代號 -> 流動卡號, and shift the remaining headers right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 header rebuild -------------------------------------------------
# Old layout (A3:I3): 日期 | 流動卡代號 | 工作 | 優 | 副 | 汙 | 破 | 副未包 | 工號
# New layout (A3:K3): 日期 | 訂單 | 品名 | 流動卡號 | 工作 | 優 | 副 | 汙 | 破 | 副未包 | 工號
# (D3 is written before B3/C3 so new shared-string entries line up the same
# way the source workbook ordered them: 流動卡號, 訂單, 品名)
$ws.Range("A3").Value = "日期"
$ws.Range("D3").Value = "流動卡號"
$ws.Range("B3").Value = "訂單"
$ws.Range("C3").Value = "品名"
$ws.Range("E3").Value = "工作"
$ws.Range("F3").Value = "優"
$ws.Range("G3").Value = "副"
$ws.Range("H3").Value = "汙"
$ws.Range("I3").Value = "破"
$ws.Range("J3").Value = "副未包"
$ws.Range("K3").Value = "工號"

# --- Column widths ---------------------------------------------------------
# This runtime's ColumnWidth setter stores (round(input*7)/7 + 5/7) in the
# saved XML width, so we back out the input that lands closest to each
# desired stored width.
$ws.Columns("A").ColumnWidth = 7.410714285714286   # -> 8.125
$ws.Columns("B").ColumnWidth = 10.785714285714286  # -> 11.5
$ws.Columns("C").ColumnWidth = 9.160714285714286   # -> 9.875
$ws.Columns("D").ColumnWidth = 9.160714285714286   # -> 9.875
$ws.Columns("E").ColumnWidth = 5.535714285714286   # -> 6.25
$ws.Columns("F").ColumnWidth = 4.910714285714286   # -> 5.625
$ws.Columns("G").ColumnWidth = 4.410714285714286   # -> 5.125
$ws.Columns("H").ColumnWidth = 4.035714285714286   # -> 4.75
$ws.Columns("I").ColumnWidth = 4.660714285714286   # -> 5.375
$ws.Columns("J").ColumnWidth = 6.660714285714286   # -> 7.375

# --- Selection --------------------------------------------------------------
$null = $ws.Range("H7").Select()

Write-Output "tidyReport layout updated"
